$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (comb1..comb3) ---
$ws.Range("B2").Value = 4
$ws.Range("E2").Value = -0.950000000000001

$ws.Range("B3").Value = 2
$ws.Range("E3").Value = -3.048814655172414

$ws.Range("B4").Value = 2
$ws.Range("E4").Value = -3.007327586206896

# --- Add new rows (comb4..comb6), copying A2's style (bold/border/center) ---
$ws.Range("A2").Copy($ws.Range("A5"))
$ws.Range("A2").Copy($ws.Range("A6"))
$ws.Range("A2").Copy($ws.Range("A7"))

$ws.Range("A5").Value = "comb4"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = -1.88

$ws.Range("A6").Value = "comb5"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = -1.8525

$ws.Range("A7").Value = "comb6"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = -1.825
